$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update URL values for the "test3 -> test4" / "test14 -> test4" environment bump
$ws.Range("M2").Value = "https://mirandakate.cabitest4.com"
$ws.Range("G2").Value = "https://test4.cliotest.com/cabicentral/control/main"
$ws.Range("J2").Value = "https://test4.cliotest.com/warehouse/control/main"
$ws.Range("A2").Value = "https://test4.cliotest.com/backoffice/control/main"

# A2's hyperlink keeps pointing at the old (test3) target, but now shows an
# explicit display/tooltip text equal to the old URL (same stale-display
# pattern already used by the G2/J2 hyperlinks).
$links = $ws.Hyperlinks
foreach ($link in $links) {
    if ($link.Range.Address() -eq '$A$2') {
        $link.TextToDisplay = "https://test3.cliotest.com/backoffice/control/main"
    }
}

# Reset the saved scroll position of the sheet view back to the top-left
# (A1) cell, clearing the stored topLeftCell="G1".
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
